# Clear original graph & create nodes
# Updates the locationName values in the "nodes" sheet to wrap the
# location name in curly/smart quotes, and moves the active selection
# to B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nodes")

$updates = @{
    "C3"  = "locationName:”Tampines”, lat: 23.23452, lon: 87.24553 "
    "C4"  = "locationName:”Pickupp HQ”, lat: 23.23452, lon: 87.24553 "
    "C5"  = "locationName:”FTZ SG”, lat: 23.23452, lon: 87.24553"
    "C6"  = "locationName:”Satsaco HQ”, lat: 23.23452, lon: 87.24553"
    "C7"  = "locationName:”Greenland HQ”, lat: 23.23452, lon: 87.24553"
    "C8"  = "locationName:”Kerry HQ”, lat: 23.23452, lon: 87.24553"
    "C9"  = "locationName:”Entrego HQ”, lat: 23.23452, lon: 87.24553"
    "C10" = "locationName:”2Go HQ”, lat: 23.23452, lon: 87.24553"
    "C11" = "locationName:”Seko HQ”, lat: 23.23452, lon: 87.24553"
    "C12" = "locationName:”GoJek HQ”, lat: 23.23452, lon: 87.24553"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws.Range("B12").Select()
